$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.931.83'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").Value = '2.237.07'
$ws.Range("E3").Value = '  +1.68%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '97.44'
$ws.Range("E5").Value = '  +16.70%  '

$ws.Range("E6").Value = '  +5.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  +1.74%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.633'
$ws.Range("E9").Value = '  +5.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.18'
$ws.Range("E10").Value = '  +8.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  +1.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.26'
$ws.Range("E12").Value = '  +14.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.30'
$ws.Range("E14").Value = '  +6.85%  '

$ws.Range("D15").Value = '2.576.70'
$ws.Range("E15").Value = '  +1.99%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.825'
$ws.Range("E16").Value = '  +5.09%  '

$ws.Range("D17").Value = '2.244.46'
$ws.Range("E17").Value = '  +3.10%  '

$ws.Range("D18").Value = '43.955.66'
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("E19").Value = '  +1.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("E20").Value = '  +4.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.87'
$ws.Range("E21").Value = '  +1.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.31'
$ws.Range("E22").Value = '  -1.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.03'
$ws.Range("E23").Value = '  +1.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("E24").Value = '  +2.89%  '

$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.40'
$ws.Range("E26").Value = '  +6.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("E27").Value = '  +12.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.53'
$ws.Range("E28").Value = '  +2.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.43'
$ws.Range("E29").Value = '  +0.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.27'
$ws.Range("E30").Value = '  +0.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.60'
$ws.Range("E31").Value = '  -0.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0910'
$ws.Range("E32").Value = '  +5.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.06'
$ws.Range("E33").Value = '  +3.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.57'
$ws.Range("E34").Value = '  +4.40%  '

$ws.Range("E35").Value = '  +1.31%  '

$ws.Range("E36").Value = '  +0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0350'
$ws.Range("E37").Value = '  -2.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.37'
$ws.Range("E38").Value = '  -3.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.60'
$ws.Range("E39").Value = '  +27.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.248'
$ws.Range("E40").Value = '  +24.39%  '

$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.18'
$ws.Range("E41").Value = '  +3.75%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.43'
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '61.96'
$ws.Range("E43").Value = '  -1.74%  '

$ws.Range("E44").Value = '  -1.04%  '

$ws.Range("E45").Value = '  +4.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.41'
$ws.Range("E46").Value = '  +1.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.28'
$ws.Range("E47").Value = '  +0.50%  '

$ws.Range("E48").Value = '  +4.06%  '

$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.435'
$ws.Range("E50").Value = '  -1.42%  '

$ws.Range("D51").Value = '2.458.15'
$ws.Range("E51").Value = '  +1.77%  '
